# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 272
$ws1.Range("F7").Value = 98
$ws1.Range("F8").Value = 263
$ws1.Range("F14").Value = 35
$ws1.Range("F19").Value = 568
$ws1.Range("F24").Value = 1977
$ws1.Range("F25").Value = 4110
$ws1.Range("F32").Value = 470
$ws1.Range("F34").Value = 122
$ws1.Range("F35").Value = 292
$ws1.Range("F36").Value = 426
$ws1.Range("F38").Value = 713
$ws1.Range("F41").Value = 422

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 272
$ws4.Range("F7").Value = 98
$ws4.Range("F8").Value = 263
$ws4.Range("F14").Value = 35
$ws4.Range("F20").Value = 568
$ws4.Range("F25").Value = 1977
$ws4.Range("F26").Value = 4110
$ws4.Range("F33").Value = 470
$ws4.Range("F35").Value = 122
$ws4.Range("F36").Value = 292
$ws4.Range("F37").Value = 426
$ws4.Range("F39").Value = 713
$ws4.Range("F42").Value = 422
